# Update column G ('K') values for Sheet1 per the regenerated save_data.
# Commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K (column G) value, as computed by the regenerated pipeline.
$kValues = @{
    2 = 0
    3 = 1
    4 = 1
    5 = 2
    6 = 2
    7 = 1
    8 = 1
    9 = 0
    11 = 2
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 3
    23 = 0
    24 = 0
    25 = 2
    26 = 2
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 3
    43 = 1
    44 = 2
    45 = 3
    46 = 1
    47 = 0
    48 = 3
    49 = 1
    50 = 1
    51 = 0
    52 = 2
    53 = 2
    54 = 1
    55 = 1
    56 = 0
    57 = 2
    58 = 1
    59 = 0
    60 = 2
    63 = 1
    64 = 2
    65 = 3
    67 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

